# Auto-generated edit script: updates cryptocurrency Price (D) and
# Volume(1h) (E) columns to match the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.894.60"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "1.650.73"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'310.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.3884"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("D8").Value = "'0.3831"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "'1.340"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'0.08438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "'23.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "'7.006"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("D15").Value = "'8.012"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "'0.00001316"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "1.651.84"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "'94.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'0.06976"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'19.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").Value = "'6.949"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'13.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "23.882.76"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "'2.444"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").Value = "'2.919"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.81%  "
$ws.Range("D27").Value = "'21.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").Value = "'153.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").Value = "'5.390"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'137.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'7.725"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").Value = "'2.486"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "1.828.41"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "'0.9919"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("D36").Value = "'0.02910"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("D37").Value = "'6.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("D38").Value = "'0.2678"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "'10.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").Value = "'0.09096"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "'0.7554"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").Value = "'13.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'16.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'0.6932"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Value = "'2.444"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").Value = "'4.095"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'0.9999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'0.08266"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").Value = "'133.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").Value = "'1.224"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.94%  "
